$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Platform Phase - add "other" option
$ws.Range("F9").Value2 = "EVT, DVT1, DVT2,PV, MP,other"

# BT Interface - add "other" option
$ws.Range("F15").Value2 = "USB,PCIe,other"

# Brand of Mouse - add "other" option
$ws.Range("F22").Value2 = "Microsoft,Logitech, Samsung, Asus,Dell,Sony,HP,other,None"

# Mouse BT - add "None" option
$ws.Range("F23").Value2 = "LE,Classic,None"

# Keyboard BT - filter build-in keyboard (set Generate UI + possible selection)
$ws.Range("E26").Value2 = "o"
$ws.Range("F26").Value2 = "LE,Classic,None"

# Brand of Keyboard - add "other" option
$ws.Range("F27").Value2 = "Microsoft,Logitech, Samsung, Asus,Dell,Sony,HP,other,None"

# Headset BT - filter build-in headset (set Generate UI + possible selection)
$ws.Range("E30").Value2 = "o"
$ws.Range("F30").Value2 = "LE,Classic,None"

# Brand of Headset - add "other" option
$ws.Range("F31").Value2 = "Microsoft,Logitech, Samsung, Asus,Dell,Sony,HP,other,None"

# Brand of Speaker - add "other" option
$ws.Range("F34").Value2 = "Microsoft,Logitech, Samsung, Asus,Dell,Sony,HP,other,None"

# Brand of Phone - add "other" option
$ws.Range("F36").Value2 = "Microsoft,Logitech, Samsung, Asus,Dell,Sony,HP,other,None"

# AC/DC mode - reorder values
$ws.Range("F63").Value2 = " DC,AC"

# Urgent Level - reorder + add "other" option
$ws.Range("F64").Value2 = "None,Fireball , P1 , P2 , P3,other"
